$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.482.80"
$ws.Range("E2").Value = "  +0.73%  "
$ws.Range("D3").Value = "1.913.73"
$ws.Range("E3").Value = "  +0.12%  "
$ws.Range("D4").Value = "'1.007"
$ws.Range("E4").Value = "  +0.51%  "
$ws.Range("D5").Value = "'325.75"
$ws.Range("E5").Value = "  +1.49%  "
$ws.Range("D6").Value = "'1.006"
$ws.Range("E6").Value = "  +0.49%  "
$ws.Range("E7").Value = "  +1.99%  "
$ws.Range("D8").Value = "'0.4071"
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("E9").Value = "  +1.39%  "
$ws.Range("D10").Value = "'1.014"
$ws.Range("E10").Value = "  +1.16%  "
$ws.Range("D11").Value = "'23.43"
$ws.Range("E11").Value = "  +4.25%  "
$ws.Range("D12").Value = "1.911.39"
$ws.Range("E12").Value = "  -2.66%  "
$ws.Range("D13").Value = "'6.016"
$ws.Range("E13").Value = "  +2.13%  "
$ws.Range("D14").Value = "'7.173"
$ws.Range("E14").Value = "  +0.55%  "
$ws.Range("D15").Value = "'90.23"
$ws.Range("E15").Value = "  +0.68%  "
$ws.Range("D16").Value = "'0.06795"
$ws.Range("E16").Value = "  +2.31%  "
$ws.Range("E17").Value = "  +0.61%  "
$ws.Range("E18").Value = "  +0.79%  "
$ws.Range("D19").Value = "'17.70"
$ws.Range("E19").Value = "  +0.17%  "
$ws.Range("E20").Value = "  +0.45%  "
$ws.Range("D21").Value = "29.509.85"
$ws.Range("D22").Value = "'5.629"
$ws.Range("E22").Value = "  +2.05%  "
$ws.Range("D23").Value = "'11.74"
$ws.Range("E23").Value = "  +2.53%  "
$ws.Range("D24").Value = "'2.180"
$ws.Range("E24").Value = "  -0.73%  "
$ws.Range("D25").Value = "2.141.12"
$ws.Range("E25").Value = "  -1.33%  "
$ws.Range("D26").Value = "'155.75"
$ws.Range("E26").Value = "  +0.47%  "
$ws.Range("D27").Value = "'6.394"
$ws.Range("E27").Value = "  +6.91%  "
$ws.Range("D28").Value = "'20.01"
$ws.Range("E28").Value = "  +1.32%  "
$ws.Range("D29").Value = "'2.110"
$ws.Range("E29").Value = "  +0.31%  "
$ws.Range("D30").Value = "'120.00"
$ws.Range("E30").Value = "  +2.11%  "
$ws.Range("D31").Value = "'1.024"
$ws.Range("E31").Value = "  -4.44%  "
$ws.Range("D32").Value = "'0.09531"
$ws.Range("E32").Value = "  +0.21%  "
$ws.Range("E33").Value = "  +2.51%  "
$ws.Range("D34").Value = "'3.563"
$ws.Range("E34").Value = "  +0.66%  "
$ws.Range("D35").Value = "'1.392"
$ws.Range("E35").Value = "  -2.08%  "
$ws.Range("D36").Value = "'0.02271"
$ws.Range("E36").Value = "  +1.11%  "
$ws.Range("D37").Value = "'0.06105"
$ws.Range("D38").Value = "'1.178"
$ws.Range("E38").Value = "  +0.17%  "
$ws.Range("D39").Value = "'0.5963"
$ws.Range("E39").Value = "  +1.90%  "
$ws.Range("D40").Value = "'10.81"
$ws.Range("E40").Value = "  +6.95%  "
$ws.Range("D41").Value = "'7.985"
$ws.Range("E41").Value = "  -3.12%  "
$ws.Range("D42").Value = "'0.1854"
$ws.Range("E42").Value = "  +1.08%  "
$ws.Range("D43").Value = "'1.279"
$ws.Range("E43").Value = "  +0.67%  "
$ws.Range("D44").Value = "'2.388"
$ws.Range("E44").Value = "  -4.38%  "
$ws.Range("E45").Value = "  +3.68%  "
$ws.Range("D46").Value = "'0.07593"
$ws.Range("E46").Value = "  -3.73%  "
$ws.Range("D47").Value = "'0.5571"
$ws.Range("E47").Value = "  +0.91%  "
$ws.Range("D48").Value = "'1.941"
$ws.Range("E48").Value = "  +1.05%  "
$ws.Range("D49").Value = "'115.88"
$ws.Range("E49").Value = "  +2.58%  "
$ws.Range("D50").Value = "'72.56"
$ws.Range("E50").Value = "  +1.98%  "
$ws.Range("E51").Value = "  +2.87%  "
